# Update cryptos list with new prices/volumes scraped on Mon Feb 20 11:36:55 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10/11, 30/31 and 35/36 swap coin identity (Coin name + Link) in addition
# to their Price/Volume updates - the coin that used to rank in the first slot
# now ranks in the second, and vice versa.
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

# Some Price cells look like plain numbers and would otherwise be silently
# re-typed as numeric by Excel's usual text -> number inference, losing the
# exact original formatting (e.g. trailing zeros such as "1.000" or "53.40").
# Force those specific cells to Text format first so the literal string is
# preserved exactly as scraped.
$textPriceCells = @("D4", "D10", "D11", "D15", "D16", "D20", "D35", "D50")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (D) and Volume(1h) (E) updates for every data row.
$ws.Range("D2").Value = "24.797.51"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.709.27"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "318.03"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "0.3921"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "1.495"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "0.9970"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").Value = "53.40"
$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("D12").Value = "0.08824"
$ws.Range("E12").Value = "  -0.71%  "

$ws.Range("D13").Value = "26.36"
$ws.Range("E13").Value = "  +11.15%  "

$ws.Range("D14").Value = "7.513"
$ws.Range("E14").Value = "  -2.85%  "

$ws.Range("D15").Value = "8.130"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "0.00001360"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").Value = "1.710.80"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "97.22"
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").Value = "0.07199"
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").Value = "20.60"
$ws.Range("E20").Value = "  +3.82%  "

$ws.Range("D21").Value = "7.318"
$ws.Range("E21").Value = "  +2.64%  "

$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("D23").Value = "14.43"
$ws.Range("E23").Value = "  -2.44%  "

$ws.Range("D24").Value = "24.796.08"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "3.021"
$ws.Range("E25").Value = "  -3.63%  "

$ws.Range("D26").Value = "2.336"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").Value = "23.15"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").Value = "168.33"
$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").Value = "5.982"
$ws.Range("E29").Value = "  +16.59%  "

$ws.Range("D30").Value = "145.91"
$ws.Range("E30").Value = "  +4.71%  "

$ws.Range("D31").Value = "8.513"
$ws.Range("E31").Value = "  -7.89%  "

$ws.Range("D32").Value = "1.899.56"
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("D33").Value = "2.217"
$ws.Range("E33").Value = "  +13.11%  "

$ws.Range("D34").Value = "0.08841"
$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("D35").Value = "0.03180"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("D36").Value = "1.057"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("D37").Value = "7.241"
$ws.Range("E37").Value = "  -8.78%  "

$ws.Range("D38").Value = "0.2837"
$ws.Range("E38").Value = "  +1.81%  "

$ws.Range("D39").Value = "0.8498"

$ws.Range("D40").Value = "10.89"
$ws.Range("E40").Value = "  -2.05%  "

$ws.Range("D41").Value = "0.09233"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").Value = "14.16"
$ws.Range("E42").Value = "  -1.99%  "

$ws.Range("D43").Value = "1.479"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("D44").Value = "17.69"
$ws.Range("E44").Value = "  +9.40%  "

$ws.Range("D45").Value = "2.705"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").Value = "0.7457"
$ws.Range("E46").Value = "  +2.84%  "

$ws.Range("D47").Value = "4.277"
$ws.Range("E47").Value = "  +1.58%  "

$ws.Range("D48").Value = "1.395"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("D49").Value = "0.9988"
$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").Value = "140.40"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "0.08294"
$ws.Range("E51").Value = "  +3.80%  "
